# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows 2-14: column C ("Förändrad") date bumped from 46065 to 46066 (+1 day)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# Rows 7,8,10,11,12,13,14: reorder entries (A, B, G columns) per updated data feed
$ws.Range("A7").Value = "A 62884-2021"
$ws.Range("B7").Value = 44504
$ws.Range("G7").Value = 0.8

$ws.Range("A8").Value = "A 25634-2025"
$ws.Range("B8").Value = 45803.59570601852
$ws.Range("G8").Value = 6

$ws.Range("A10").Value = "A 25015-2023"
$ws.Range("B10").Value = 45085.6989699074
$ws.Range("G10").Value = 1.8

$ws.Range("A11").Value = "A 19922-2025"
$ws.Range("B11").Value = 45771.63034722222
$ws.Range("G11").Value = 10.1

$ws.Range("A12").Value = "A 60024-2025"
$ws.Range("B12").Value = 45992
$ws.Range("G12").Value = 1.1

$ws.Range("A13").Value = "A 3402-2026"
$ws.Range("B13").Value = 46042.39047453704
$ws.Range("G13").Value = 5.5

$ws.Range("A14").Value = "A 14271-2021"
$ws.Range("B14").Value = 44278
$ws.Range("G14").Value = 6.7
